# Update odds values on Sheet1 to reflect the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "Q3"  = 2.05
    "R3"  = 1.8

    "G6"  = 1.36
    "H6"  = 5
    "M6"  = 1.06
    "N6"  = 10
    "AE6" = 9
    "AJ6" = 8.5

    "S7"  = 1.62
    "T7"  = 2.3

    "O9"  = 1.18
    "P9"  = 4.5
    "S9"  = 1.62
    "T9"  = 2.25
    "U9"  = 2
    "V9"  = 1.8

    "G10"  = 2.8
    "I10"  = 2.6
    "J10"  = 3.5
    "L10"  = 3.25
    "T10"  = 1.7
    "AC10" = 8.5
    "AE10" = 11
    "AF10" = 29

    "G11"  = 1.45
    "H11"  = 5
    "I11"  = 6
    "J11"  = 1.91
    "L11"  = 6
    "S11"  = 1.5
    "AA11" = 1.7
    "AB11" = 2.05
    "AC11" = 9
    "AJ11" = 9.5
    "AN11" = 34
    "AP11" = 67
    "AR11" = 41

    "M12"  = 1.04
    "N12"  = 10
    "O12"  = 1.27
    "T12"  = 1.75
    "X12"  = 1.25
    "AF12" = 41
    "AI12" = 9
    "AQ12" = 19
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
